$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3143
$ws1.Range("F3").Value = 536
$ws1.Range("F4").Value = 1101
$ws1.Range("F5").Value = 89
$ws1.Range("F6").Value = 39
$ws1.Range("F8").Value = 40
$ws1.Range("F9").Value = 1132
$ws1.Range("F10").Value = 15817
$ws1.Range("F11").Value = 248
$ws1.Range("F14").Value = 6207
$ws1.Range("F16").Value = 110
$ws1.Range("F17").Value = 68
$ws1.Range("F20").Value = 1264
$ws1.Range("F21").Value = 30
$ws1.Range("F23").Value = 21
$ws1.Range("F25").Value = 4
$ws1.Range("F27").Value = 870
$ws1.Range("F28").Value = 33
$ws1.Range("F29").Value = 5008
$ws1.Range("F30").Value = 491
$ws1.Range("F31").Value = 11102
$ws1.Range("F32").Value = 1234
$ws1.Range("F33").Value = 15
$ws1.Range("F34").Value = 129
$ws1.Range("F35").Value = 181
$ws1.Range("F37").Value = 267

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3143
$ws4.Range("F4").Value = 536
$ws4.Range("F5").Value = 1101
$ws4.Range("F6").Value = 89
$ws4.Range("F7").Value = 39
$ws4.Range("F9").Value = 40
$ws4.Range("F10").Value = 1132
$ws4.Range("F11").Value = 15817
$ws4.Range("F12").Value = 248
$ws4.Range("F15").Value = 6207
$ws4.Range("F17").Value = 110
$ws4.Range("F18").Value = 68
$ws4.Range("F21").Value = 1264
$ws4.Range("F22").Value = 30
$ws4.Range("F24").Value = 21
$ws4.Range("F26").Value = 4
$ws4.Range("F28").Value = 870
$ws4.Range("F29").Value = 33
$ws4.Range("F30").Value = 5008
$ws4.Range("F31").Value = 491
$ws4.Range("F33").Value = 11102
$ws4.Range("F34").Value = 1234
$ws4.Range("F35").Value = 15
$ws4.Range("F36").Value = 129
$ws4.Range("F37").Value = 181
$ws4.Range("F39").Value = 267
